$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2508.7273
$ws.Range("I40").Value = 1932.6666
$ws.Range("J40").Value = 3200
$ws.Range("K40").Value = 1932.6666
$ws.Range("L40").Value = 3200
$ws.Range("M40").Value = -1757.6666
$ws.Range("N40").Value = -3550

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 251.11111
$ws.Range("I53").Value = 241.16667
$ws.Range("K53").Value = 241.16667
$ws.Range("M53").Value = 395.83333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 110964.29
$ws.Range("I70").Value = 1159.3334
$ws.Range("K70").Value = 3478.0002
$ws.Range("M70").Value = -3208.0002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 110964.29
$ws.Range("I73").Value = 1159.3334
$ws.Range("K73").Value = 3478.0002
$ws.Range("M73").Value = -2542.0002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1324.5333
$ws.Range("I92").Value = 1214.25
$ws.Range("J92").Value = 1450.5714
$ws.Range("K92").Value = 1214.25
$ws.Range("L92").Value = 1450.5714
$ws.Range("M92").Value = 33.75
$ws.Range("N92").Value = -3946.5714

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 880.9545000000001
$ws.Range("I98").Value = 789.05
$ws.Range("K98").Value = 789.05
$ws.Range("M98").Value = 708.95

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 880.9545000000001
$ws.Range("I122").Value = 789.05
$ws.Range("K122").Value = 2367.15
$ws.Range("M122").Value = 82.85000000000036

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1154.3143
$ws.Range("I132").Value = 1120.2333
$ws.Range("J132").Value = 1358.8
$ws.Range("K132").Value = 3360.699900000001
$ws.Range("L132").Value = 4076.4
$ws.Range("M132").Value = -830.6999000000005
$ws.Range("N132").Value = -9136.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1777.6666
$ws.Range("I135").Value = 1083.3334
$ws.Range("J135").Value = 3166.3333
$ws.Range("K135").Value = 9750.000599999999
$ws.Range("L135").Value = 28496.9997
$ws.Range("M135").Value = -7215.000599999999
$ws.Range("N135").Value = -33566.9997

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4354
$ws.Range("J138").Value = 4751.582
$ws.Range("L138").Value = 14254.746
$ws.Range("N138").Value = -24534.746

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 878.5454999999999
$ws.Range("I2").Value = 841.75
$ws.Range("K2").Value = 841.75
$ws.Range("M2").Value = -728.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15830.44
$ws.Range("I32").Value = 7001.4243
$ws.Range("K32").Value = 7001.4243
$ws.Range("M32").Value = -6714.4243

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5198.8
$ws.Range("I45").Value = 1874
$ws.Range("K45").Value = 1874
$ws.Range("M45").Value = -1497

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2954.2
$ws.Range("I61").Value = 2942.875
$ws.Range("J61").Value = 2999.5
$ws.Range("K61").Value = 2942.875
$ws.Range("L61").Value = 2999.5
$ws.Range("M61").Value = -2730.875
$ws.Range("N61").Value = -3423.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5090.9443
$ws.Range("I63").Value = 3911.625
$ws.Range("K63").Value = 3911.625
$ws.Range("M63").Value = -3225.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 5090.9443
$ws.Range("I66").Value = 3911.625
$ws.Range("K66").Value = 19558.125
$ws.Range("M66").Value = -16126.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4132.727
$ws.Range("I74").Value = 1470.6666
$ws.Range("K74").Value = 1470.6666
$ws.Range("M74").Value = -596.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4132.727
$ws.Range("I77").Value = 1470.6666
$ws.Range("K77").Value = 7353.333000000001
$ws.Range("M77").Value = -2985.333000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 4649.75
$ws.Range("I110").Value = 5178.2856
$ws.Range("J110").Value = 950
$ws.Range("K110").Value = 5178.2856
$ws.Range("L110").Value = 950
$ws.Range("M110").Value = -3133.2856
$ws.Range("N110").Value = -5040

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 878.5454999999999
$ws.Range("I116").Value = 841.75
$ws.Range("K116").Value = 841.75
$ws.Range("M116").Value = 1452.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1775.0625
$ws.Range("I132").Value = 1683.1724
$ws.Range("J132").Value = 2663.3333
$ws.Range("K132").Value = 5049.5172
$ws.Range("L132").Value = 7989.999899999999
$ws.Range("M132").Value = -2519.5172
$ws.Range("N132").Value = -13049.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2954.2
$ws.Range("I136").Value = 2942.875
$ws.Range("J136").Value = 2999.5
$ws.Range("K136").Value = 8828.625
$ws.Range("L136").Value = 8998.5
$ws.Range("M136").Value = -6278.625
$ws.Range("N136").Value = -14098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 878.5454999999999
$ws.Range("I3").Value = 841.75
$ws.Range("K3").Value = 841.75
$ws.Range("M3").Value = -727.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4503.5
$ws.Range("I134").Value = 4503.5
$ws.Range("K134").Value = 13510.5
$ws.Range("M134").Value = -10975.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3716.375
$ws.Range("I58").Value = 2136.875
$ws.Range("K58").Value = 2136.875
$ws.Range("M58").Value = -1933.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 292.33334
$ws.Range("I107").Value = 223.92857
$ws.Range("K107").Value = 223.92857
$ws.Range("M107").Value = 1696.07143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 72874.75
$ws.Range("J125").Value = 72874.75
$ws.Range("L125").Value = 72874.75
$ws.Range("N125").Value = -77794.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3781.2778
$ws.Range("I132").Value = 2951.0908
$ws.Range("K132").Value = 8853.2724
$ws.Range("M132").Value = -6323.2724

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 92238.5
$ws.Range("J133").Value = 91990
$ws.Range("L133").Value = 91990
$ws.Range("N133").Value = -97050

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4818.1665
$ws.Range("I134").Value = 3448.3333
$ws.Range("K134").Value = 10344.9999
$ws.Range("M134").Value = -7809.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3716.375
$ws.Range("I136").Value = 2136.875
$ws.Range("K136").Value = 6410.625
$ws.Range("M136").Value = -3860.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 55556656
$ws.Range("I64").Value = 833
$ws.Range("J64").Value = 83334570
$ws.Range("K64").Value = 2499
$ws.Range("L64").Value = 250003710
$ws.Range("M64").Value = -2229
$ws.Range("N64").Value = -250004250

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 55556656
$ws.Range("I67").Value = 833
$ws.Range("J67").Value = 83334570
$ws.Range("K67").Value = 2499
$ws.Range("L67").Value = 250003710
$ws.Range("M67").Value = -1563
$ws.Range("N67").Value = -250005582

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3295.158
$ws.Range("J131").Value = 3555.4
$ws.Range("L131").Value = 10666.2
$ws.Range("N131").Value = -20746.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 55000
$ws.Range("I64").Value = 55000
$ws.Range("K64").Value = 55000
$ws.Range("M64").Value = -54752

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H67").Value = 55000
$ws.Range("I67").Value = 55000
$ws.Range("K67").Value = 55000
$ws.Range("M67").Value = -54142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1608.1786
$ws.Range("I97").Value = 1777
$ws.Range("J97").Value = 1186.125
$ws.Range("K97").Value = 1777
$ws.Range("L97").Value = 1186.125
$ws.Range("M97").Value = -1281
$ws.Range("N97").Value = -2178.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 25338.357
$ws.Range("J136").Value = 25338.357
$ws.Range("L136").Value = 76015.071
$ws.Range("N136").Value = -81115.071

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 64324.75
$ws.Range("J141").Value = 64324.75
$ws.Range("L141").Value = 64324.75
$ws.Range("N141").Value = -74684.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2874.5
$ws.Range("J68").Value = 2874.5
$ws.Range("L68").Value = 2874.5
$ws.Range("N68").Value = -4372.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2874.5
$ws.Range("J71").Value = 2874.5
$ws.Range("L71").Value = 14372.5
$ws.Range("N71").Value = -21860.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 633.3333
$ws.Range("I93").Value = 500
$ws.Range("J93").Value = 900
$ws.Range("K93").Value = 500
$ws.Range("L93").Value = 900
$ws.Range("M93").Value = 748
$ws.Range("N93").Value = -3396

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5314.0303
$ws.Range("I132").Value = 4916.773
$ws.Range("J132").Value = 6108.5454
$ws.Range("K132").Value = 14750.319
$ws.Range("L132").Value = 18325.6362
$ws.Range("M132").Value = -12220.319
$ws.Range("N132").Value = -23385.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3852.75
$ws.Range("I136").Value = 3673.3
$ws.Range("K136").Value = 11019.9
$ws.Range("M136").Value = -8469.900000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1544.4445
$ws.Range("I100").Value = 1362.5
$ws.Range("K100").Value = 2725
$ws.Range("M100").Value = -2184

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1088.1875
$ws.Range("I107").Value = 416.23077
$ws.Range("K107").Value = 1248.69231
$ws.Range("M107").Value = 671.3076900000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1148.7273
$ws.Range("I113").Value = 910
$ws.Range("J113").Value = 1472.7142
$ws.Range("K113").Value = 2730
$ws.Range("L113").Value = 4418.142599999999
$ws.Range("M113").Value = -560
$ws.Range("N113").Value = -8758.142599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1843.4
$ws.Range("I132").Value = 1866.75
$ws.Range("J132").Value = 1750
$ws.Range("K132").Value = 5600.25
$ws.Range("L132").Value = 5250
$ws.Range("M132").Value = -3070.25
$ws.Range("N132").Value = -10310

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 57303.445
$ws.Range("I136").Value = 1264.5454
$ws.Range("J136").Value = 145364.58
$ws.Range("K136").Value = 3793.6362
$ws.Range("L136").Value = 436093.74
$ws.Range("M136").Value = -1243.6362
$ws.Range("N136").Value = -441193.74
